$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-CellText 'A4' '7abaCsl-b5MC'
Set-CellText 'B4' 'Tolkien'
Set-CellText 'C4' 'Pocos escritores vivos han conseguido atraer tanta atencion e interes, despertar tanta pasion y entusiasmo, suscitar tanta controversia y adulacion como J.R.R Tolkien. Curiosamente, nadie ha explorado los origenes de la obra cumbre de Tolkien, El Senor de los Anillos, ni ha intentado relacionarla con la larga saga de historias epicas en las que se inspira y de las que es, a la vez, ejemplo.El origen de El Senor de los Anillos es un ameno recorrido por las sagas heroicas, desde los mitos griegos hasta las mitologias nordicas, y los grandes relatos que, a traves del tiempo, se han convertido en la fuente e inspiracion de uno de los mayores narradores del siglo XX: J.R.R. Tolkien.El libro de Lin Carter constituye una lectura absorbente para los millones de seguidores de El Senor de los anillos y ofrece una fascinante introduccion para quienes deseen aproximarse por primera vez al universo magico de Tolkien.'
Set-CellText 'D4' '2002'
Set-CellText 'E4' 'J. R. R. Tolkien, Lin Carter'
Set-CellText 'A6' 'GAcCAAAACAAJ'
Set-CellText 'B6' 'Los mundos mágicos de El señor de los anillos'
Set-CellText 'C6' 'Desconocido'
Set-CellText 'D6' '2003'
Set-CellText 'E6' 'David Colbert'
Set-CellText 'A8' 'IOyjMQEACAAJ'
Set-CellText 'B8' 'The Power of the Ring'
Set-CellText 'C8' 'Digging deep into J. R. R. Tolkien''s spiritual biography--his religious scholarship and his love of both Christian and pagan myth--Stratford Caldecott offers a critical study of how the acclaimed author effectively created a vivid Middle Earth using the familiar rites and ceremonies of human history. And while readers and moviegoers alike may appreciate the fantasy world of The Hobbit and the Lord of the Rings trilogy, few know that in life, Tolkien was a devout Roman Catholic and that the characters, the events, and the general morality of each novel are informed by the dogmas of his faith. Revised and updated, this acclaimed study of Tolkien''s achievement--previously released as Secret Fire in the UK--includes commentary on Peter Jackson''s film adaptations and explores many of the fascinating stories and letters published after Tolkien''s death.'
Set-CellText 'D8' '2012-12-04'
Set-CellText 'E8' 'Stratford Caldecott'
Set-CellText 'A9' 'JUdOAAAACAAJ'
Set-CellText 'B9' 'The Lord of the Rings Sketchbook'
Set-CellText 'C9' '"In The Lord of the Rings Sketchbook Alan Lee reveals in pictures and in words how he created the watercolor paintings for the special centenary edition of The Lord of the Rings. These images would prove so powerful and evocative that they would eventually define the look of Peter Jackson''s movie trilogy and earn him a coveted Academy Award." "The book is filled with more than 150 of his sketches and early conceptual pieces showing how the project progressed from idea to finished art. It also contains a selection of full-page paintings reproduced in full color, together with numerous examples of previously unseen conceptual art produced for the films and many new works drawn specially for this book." "The Lord of the Rings Sketchbook provides an insight into the imagination of the man who painted Tolkien''s vision, first on the page and then in three dimensions on the movie screen. It will also be of interest to many of the thousands of people who have bought the illustrated Lord of the Rings as well as to budding artists who want to unlock the secrets of book illustration."--BOOK JACKET.'
Set-CellText 'D9' '2005'
Set-CellText 'E9' 'Alan Lee'
Set-CellText 'A10' 'LCZvVRqH-m8C'
Set-CellText 'B10' 'El señor de los anillos'
Set-CellText 'C10' 'Desconocido'
Set-CellText 'D10' '1999-05'
Set-CellText 'E10' 'Terry Donaldson'
Set-CellText 'A11' 'LvsQ34A1fOMC'
Set-CellText 'B11' 'El Señor de los Anillos no 03/03 El Retorno del Rey (edición revisada)'
Set-CellText 'C11' 'La tercera entrega de la trilogía El Señor de los Anillos. Los ejércitos del Señor Oscuro van extendiendo cada vez más su maléfica sombra por la Tierra Media. Hombres, elfos y enanos unen sus fuerzas para presentar batalla a Sauron y sus huestes. Ajenos a estos preparativos, Frodo y Sam siguen adentrándose en el país de Mordor en su heroico viaje para destruir el Anillo de Poder en las Grietas del Destino. «Un final triunfante... un gran trabajo, tanto en la concepción como en la ejecución.» —Daily Telegraph «Un trabajo extraordinariamente imaginativo, parte saga, parte alegoría, y emocionante en su totalidad.» —The Times'
Set-CellText 'D11' '2010-07-15'
Set-CellText 'E11' 'J. R. R. Tolkien'
Set-CellText 'A12' 'Ndgf0AEACAAJ'
Set-CellText 'B12' 'El señor de los anillos'
Set-CellText 'C12' 'Desconocido'
Set-CellText 'D12' '2007'
Set-CellText 'E12' 'John Ronald Reuel Tolkien'
Set-CellText 'A13' 'PC8AAAAACAAJ'
Set-CellText 'B13' 'El Señor de los Anillos'
Set-CellText 'C13' 'Ediciones Minotauro se complace en presentar la publicacion oficial y autorizada de La Comunidad del Anillo. Clbum de la pelicula, coincidiendo con el estreno de la primera parte de la trilogia epica de Peter Jackson.'
Set-CellText 'D13' '2001'
Set-CellText 'E13' 'Jude Fisher'
Set-CellText 'A18' 'UfYGAAAACAAJ'
Set-CellText 'B18' 'El señor de los anillos'
Set-CellText 'D18' '2002-02'
Set-CellText 'E18' 'Kurt D. Bruner, Jim Ware'
Set-CellText 'A19' 'WmdWtQAACAAJ'
Set-CellText 'B19' 'El Señor de los anillos'
Set-CellText 'D19' '2002'
Set-CellText 'A20' 'ZVwX0QEACAAJ'
Set-CellText 'B20' 'El Señor de los Anillos'
Set-CellText 'C20' 'Desconocido'
Set-CellText 'D20' '1985'
Set-CellText 'A25' 'gHzeCYjKGBsC'
Set-CellText 'B25' 'El señor de los anillos y la filosofía'
Set-CellText 'C25' '¿Nos resultaría aburrida la vida eterna en caso de que pudiéramos alcanzarla? ¿Qué virtudes son necesarias para que el poder no corrompa a quien lo tiene? ¿Se revelará la naturaleza contra la tecnología? Si un ente cae en el bosque y nadie lo oye, ¿hace ruido? Gregory Bassham y Eric Bronson dan respuesta a éstas y otras cuestiones planteadas en El Señor de los Anillos, la epopeya de J. R. R. Tolkien que, desde su publicación, ha vendido más de cincuenta millones de ejemplares y cuya riqueza y complejidad ha tejido una red cada vez más amplia de seguidores. Si la intención de Tolkien era crear «una historia que mantendría la atención del lector [...] y que a veces quizá lo excitaría o lo conmovería profundamente», consiguió también trasladar a su obra algunas de sus propias inquietudes filosóficas: la lucha del bien contra el mal, la oposición entre destino y libre albedrío, la búsqueda de la felicidad o la vida después de la muerte. Esta colección de ensayos además de ofrecernos una comprensión más fiel de las cuestiones que inspiran y nutren El Señor de los Anillos, supone una amena introducción a la filosofía, la religión y la mitología y una excelente carta de presentación a la cosmovisión de pensadores como Platón, Aristóteles o Nietzsche. Porque, como dice el mago filósofo Gandalf: "Si has estado estos días con las orejas tapadas y la mente dormida, ¡es hora de que despiertes!'
Set-CellText 'D25' '2010-05-20'
Set-CellText 'E25' 'Gregory Bassham, Eric Bronson'
Set-CellText 'A26' 'jZgjyAEACAAJ'
Set-CellText 'B26' 'The Return of the King'
Set-CellText 'C26' 'Desconocido'
Set-CellText 'D26' '1967'
Set-CellText 'A27' 'ld5GswEACAAJ'
Set-CellText 'B27' 'The Fellowship of the Ring'
Set-CellText 'C27' 'Frodo the hobbit and a band of warriors from the different kingdoms set out to destroy the Ring of Power before the evil Sauron grasps control.'
Set-CellText 'D27' '2005'
Set-CellText 'A28' 'neSkMQEACAAJ'
Set-CellText 'B28' 'El Señor de Los Anillos 1.'
Set-CellText 'C28' 'La Compania se ha disuelto y sus integrantes emprenden caminos separados. Frodo y Sam continuan solos su viaje a lo largo del rio Anduin, perseguidos por la sombra misteriosa de un ser extrano que tambien ambiciona la posesion del Anillo. Mientras, hombres, elfos y enanos se preparan para la batalla final contra las fuerzas del Senor del Mal.'
Set-CellText 'D28' '2012-11-13'
Set-CellText 'E28' 'John Ronald Reuel Tolkien'
Set-CellText 'A29' 'o5WfPwAACAAJ'
Set-CellText 'B29' 'El Señor de los anillos'
Set-CellText 'C29' 'Desconocido'
Set-CellText 'A30' 'q0JyPwAACAAJ'
Set-CellText 'B30' 'El Señor de Los Anillos, I'
Set-CellText 'C30' 'En la adormecida e idílica Comarca, un joven hobbit recibe un encargo: custodiar el Anillo Único y emprender el viaje para su destrucción en las Grietas del Destino. Acompañado por magos, hombres, elfos y enanos, atravesará la Tierra Media y se internará en las sombras de Mordor, perseguido siempre por las huestes de Sauron, el Señor Oscuro, dispuesto a recuperar su creación para establecer el dominio definitivo del Mal.'
Set-CellText 'D30' '2002'
Set-CellText 'E30' 'J. R. R. Tolkien'
